$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New data rows (form submissions) ---
# Existing rows 2-4 already contain their data; rows 5-10 are new.
$rows = @{
  5  = @('1898367659','2024-12-01 10:23:25','Нет','Нет')
  6  = @('1898367632','2024-12-01 10:23:29','Нет','Нет')
  7  = @('1898367594','2024-12-10 15:00:00','Да','Да')
  8  = @('1898367659','2024-12-10 15:00:01','Нет','Нет')
  9  = @('1898367660','2024-12-15 23:02:02','Да','Да')
  10 = @('1898367665','2024-12-15 23:02:03','Нет','Нет')
}

# Apply text number format ("@", numFmtId 49) to all data rows (2-10) first,
# so the cells (including the new, not-yet-populated ones) get the proper
# style before values are written - this mirrors how the workbook stores
# the ID / timestamp columns as text rather than numbers / dates.
$ws.Range("A2:D10").NumberFormat = "@"

# Fill in the new rows
foreach ($r in 5..10) {
  $vals = $rows[$r]
  $ws.Cells.Item($r,1).Value = $vals[0]
  $ws.Cells.Item($r,2).Value = $vals[1]
  $ws.Cells.Item($r,3).Value = $vals[2]
  $ws.Cells.Item($r,4).Value = $vals[3]
}

# --- Column widths (approximate character widths used by the author) ---
$ws.Columns.Item(1).ColumnWidth = 22.5703125
$ws.Columns.Item(2).ColumnWidth = 21.140625
$ws.Columns.Item(3).ColumnWidth = 32.42578125
$ws.Columns.Item(4).ColumnWidth = 33.85546875

# --- Selection moved to C12 ---
$ws.Range("C12").Select() | Out-Null

# --- Workbook window position (best effort; engine may not persist this) ---
$wb.Windows.Item(1).Left = 3900
$wb.Windows.Item(1).Top = 1695

Write-Host "Edit applied"
